# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) for the 153d35a8 source-file row
# (row 2) on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 16:37:38"
$wsZhCn.Range("H2").Value = "2016-03-17 16:37:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 16:37:42"
$wsDeDe.Range("H2").Value = "2016-03-17 16:38:06"
